# Apply the "LOOKUP.xlsx" edit:
#  - Move the Neck Size / Jacket Size lookup table on Sheet1 from columns E:F to H:I
#  - Add VLOOKUP formulas in Sheet1!C4:C17 that look up the jacket size from the table
#  - Add VLOOKUP formulas in Sheet2!E4:E15 that look up the bulk-buy discount from the table
#  - Repoint the workbook-level defined names (Large, Largest, Medium, Small,
#    Threshold_1..3) from column F to column I on Sheet1
#  - Update selections on both sheets

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet1: move the lookup table from E4:F11 to H4:I11
# ---------------------------------------------------------------------------
$ws1.Range("E4:F11").Copy($ws1.Range("H4"))
$ws1.Range("E4:F11").Clear()

# Give the moved columns (H, I) roughly the same display width the table had
# when it lived in E:F.
$ws1.Columns.Item(8).ColumnWidth = 8
$ws1.Columns.Item(9).ColumnWidth = 8.8333333333333339

# ---------------------------------------------------------------------------
# Sheet1: add the VLOOKUP formulas that report a dog's jacket size
# ---------------------------------------------------------------------------
$ws1.Range("C4").Formula = "=VLOOKUP(B4,`$H`$5:`$I`$11,2)"
$ws1.Range("C5:C17").Formula = "=VLOOKUP(B5,`$H`$5:`$I`$11,2)"

# ---------------------------------------------------------------------------
# Sheet2: add the VLOOKUP formulas that report the bulk-buy discount
# ---------------------------------------------------------------------------
$ws2.Range("E4").Formula = "=VLOOKUP(C4,`$H`$5:`$I`$11,2)"
$ws2.Range("E5:E15").Formula = "=VLOOKUP(C5,`$H`$5:`$I`$11,2)"

# ---------------------------------------------------------------------------
# Workbook: repoint defined names from Sheet1!$F$n to Sheet1!$I$n
# ---------------------------------------------------------------------------
$wb.Names.Item("Large").RefersTo = "=Sheet1!`$I`$5"
$wb.Names.Item("Largest").RefersTo = "=Sheet1!`$I`$4"
$wb.Names.Item("Medium").RefersTo = "=Sheet1!`$I`$6"
$wb.Names.Item("Small").RefersTo = "=Sheet1!`$I`$7"
$wb.Names.Item("Threshold_1").RefersTo = "=Sheet1!`$I`$9"
$wb.Names.Item("Threshold_2").RefersTo = "=Sheet1!`$I`$10"
$wb.Names.Item("Threshold_3").RefersTo = "=Sheet1!`$I`$11"

# ---------------------------------------------------------------------------
# Selections: Sheet1 selects C4:C17, Sheet2 selects B22 and stays the
# selected/active tab (matching the workbook's saved view state).
# ---------------------------------------------------------------------------
$ws1.Range("C4:C17").Select() | Out-Null
$ws2.Range("B22").Select() | Out-Null
